$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to a cleaner name
$ws.Name = "Product_List"

# Insert a new "Serial Number" column right after Product ID (new column B),
# pushing Product Name/Warehouse #/Quantity/Selling Price/Cost Price one
# column to the right.
$ws.Columns.Item(2).Insert()

# Header for the newly inserted column
$ws.Cells.Item(1, 2).Value = "Serial Number"

# Serial numbers for each of the 6 product rows
$ws.Cells.Item(2, 2).Value = "SNY110"
$ws.Cells.Item(3, 2).Value = "LG220"
$ws.Cells.Item(4, 2).Value = "SMSUNG330"
$ws.Cells.Item(5, 2).Value = "SNY440"
$ws.Cells.Item(6, 2).Value = "BOSSSS556"
$ws.Cells.Item(7, 2).Value = "PONEER665"

# Clarify the price headers (now columns F and G after the insert) by
# calling out that the values are in dollars
$ws.Cells.Item(1, 6).Value = "Selling Price ($)"
$ws.Cells.Item(1, 7).Value = "Cost Price ($)"

# Drop the old "$#,##0" custom currency formatting on the price columns
# (the header already states the unit, so plain numbers are used instead)
$ws.Range("F2:G7").NumberFormat = "General"

# Resize every used column so the data is aligned and easy to read
$ws.Cells.EntireColumn.AutoFit()

# Restore the cursor/selection position recorded in the saved workbook
$ws.Range("I11").Select() | Out-Null
